$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $ref, $val)
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "56.399.50"
Set-TextCell $ws "E2" "  +5.23%  "

Set-TextCell $ws "D3" "2.482.36"
Set-TextCell $ws "E3" "  +3.53%  "

Set-TextCell $ws "E4" "  +0.13%  "

Set-TextCell $ws "D5" "487.98"
Set-TextCell $ws "E5" "  +5.66%  "

Set-TextCell $ws "D6" "147.35"
Set-TextCell $ws "E6" "  +12.80%  "

Set-TextCell $ws "D7" "0.996"
Set-TextCell $ws "E7" "  -0.10%  "

Set-TextCell $ws "D8" "0.511"
Set-TextCell $ws "E8" "  +5.21%  "

Set-TextCell $ws "D9" "2.492.16"
Set-TextCell $ws "E9" "  +3.56%  "

Set-TextCell $ws "D10" "5.78"
Set-TextCell $ws "E10" "  +9.41%  "

Set-TextCell $ws "D11" "0.0969"
Set-TextCell $ws "E11" "  +2.85%  "

Set-TextCell $ws "D12" "0.332"
Set-TextCell $ws "E12" "  +6.43%  "

Set-TextCell $ws "E13" "  +1.66%  "

Set-TextCell $ws "D14" "2.914.15"
Set-TextCell $ws "E14" "  +2.87%  "

Set-TextCell $ws "D15" "56.375.37"
Set-TextCell $ws "E15" "  +5.52%  "

Set-TextCell $ws "D16" "21.07"
Set-TextCell $ws "E16" "  +7.91%  "

Set-TextCell $ws "D17" "0.0000136"
Set-TextCell $ws "E17" "  +3.07%  "

Set-TextCell $ws "D18" "2.486.83"
Set-TextCell $ws "E18" "  +2.46%  "

Set-TextCell $ws "E19" "  +8.79%  "

Set-TextCell $ws "D20" "10.08"
Set-TextCell $ws "E20" "  +7.92%  "

Set-TextCell $ws "D21" "317.97"
Set-TextCell $ws "E21" "  +3.92%  "

Set-TextCell $ws "D22" "0.998"
Set-TextCell $ws "E22" "  +0.32%  "

Set-TextCell $ws "D23" "5.81"
Set-TextCell $ws "E23" "  +9.34%  "

Set-TextCell $ws "D24" "58.49"
Set-TextCell $ws "E24" "  +4.97%  "

Set-TextCell $ws "D25" "0.411"
Set-TextCell $ws "E25" "  +7.24%  "

Set-TextCell $ws "D27" "0.163"
Set-TextCell $ws "E27" "  +6.89%  "

Set-TextCell $ws "D28" "2.585.60"
Set-TextCell $ws "E28" "  +3.01%  "

Set-TextCell $ws "D29" "7.65"
Set-TextCell $ws "E29" "  +8.54%  "

Set-TextCell $ws "D30" "0.0₃0794"
Set-TextCell $ws "E30" "  +11.56%  "

Set-TextCell $ws "E31" "  -0.19%  "

Set-TextCell $ws "E32" "  +2.70%  "

Set-TextCell $ws "D33" "18.19"
Set-TextCell $ws "E33" "  +3.49%  "

Set-TextCell $ws "D34" "1.51"
Set-TextCell $ws "E34" "  +6.04%  "

Set-TextCell $ws "D35" "5.20"
Set-TextCell $ws "E35" "  +4.97%  "

Set-TextCell $ws "E36" "  +9.09%  "

Set-TextCell $ws "D37" "3.75"
Set-TextCell $ws "E37" "  +7.30%  "

Set-TextCell $ws "D38" "0.863"
Set-TextCell $ws "E38" "  +8.54%  "

Set-TextCell $ws "D39" "34.19"
Set-TextCell $ws "E39" "  +4.90%  "

Set-TextCell $ws "D40" "3.51"
Set-TextCell $ws "E40" "  +8.61%  "

Set-TextCell $ws "E41" "  +7.04%  "

Set-TextCell $ws "B42" "Mantle"
Set-TextCell $ws "C42" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D42" "0.611"
Set-TextCell $ws "E42" "  +3.87%  "

Set-TextCell $ws "B43" "FirstDigitalUSD"
Set-TextCell $ws "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D43" "0.995"
Set-TextCell $ws "E43" "  +0.53%  "

Set-TextCell $ws "E44" "  +8.39%  "

Set-TextCell $ws "D45" "4.78"
Set-TextCell $ws "E45" "  +15.37%  "

Set-TextCell $ws "D46" "0.0927"
Set-TextCell $ws "E46" "  +6.98%  "

Set-TextCell $ws "D47" "258.47"
Set-TextCell $ws "E47" "  +15.49%  "

Set-TextCell $ws "E48" "  +5.98%  "

Set-TextCell $ws "E49" "  +0.59%  "

Set-TextCell $ws "D50" "17.59"
Set-TextCell $ws "E50" "  +7.71%  "

Set-TextCell $ws "D51" "1.871.92"
Set-TextCell $ws "E51" "  -2.86%  "
